# Update Leaves and delete module notification
# Append two new form-response rows to the "Responses" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 -----------------------------------------------------------
$ws.Range("A6").Value = "05/04/2024 08:48:16"
$ws.Range("B6").Value = "trabelsi.ahmed.1@esprit.tn"
$ws.Range("C6").Value = "ahmed"
$ws.Range("D6").Value = "ca va"
$ws.Range("E6").Value = "oui"
$ws.Range("F6").Value = "oui"
$ws.Range("G6").Value = "oui"

# "10000" must be stored as text (matches the other plain-text survey
# answers), so force a text quote-prefix then drop the formatting again
# so the cell keeps the workbook's default style.
$ws.Range("H6").Value = "'10000"
$ws.Range("H6").ClearFormats()

$ws.Range("I6").Value = "non"
$ws.Range("J6").Value = "oui"

# --- Row 7 -----------------------------------------------------------
$ws.Range("A7").Value = "23/04/2024 22:11:13"
$ws.Range("B7").Value = "aymen.nefzi@esprit.tn"
$ws.Range("C7").Value = "aaa"
$ws.Range("D7").Value = "aaaa"
$ws.Range("E7").Value = "aaaa"
$ws.Range("F7").Value = "aaaa"
$ws.Range("G7").Value = "aaaa"
$ws.Range("H7").Value = "aaaa"
$ws.Range("I7").Value = "aaaa"
$ws.Range("J7").Value = "aaaa"
